$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header C1: "result_clean" -> "result"
$ws.Range("C1").Value = "result"

# Capitalize "abgelehnt" -> "Abgelehnt" for rows 2..21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "Abgelehnt"
}

# Capitalize "angenommen" -> "Angenommen" for rows 22..41
for ($r = 22; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = "Angenommen"
}
